$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B for "Country", pushing existing Budget column to C
$ws.Range("B1").EntireColumn.Insert()

# New header
$ws.Range("B1").Value = "Country"

# Seed the shared-string table with the country codes in the same order they
# were first introduced in the source document (us, br, ca, aunz) using a
# scratch area, then drop the values into their real positions below. This
# keeps the shared-string indices stable without depending on row-write order.
$ws.Range("Z1").Value = "us"
$ws.Range("Z2").Value = "br"
$ws.Range("Z3").Value = "ca"
$ws.Range("Z4").Value = "aunz"

$countryUs = $ws.Range("Z1").Value2
$countryBr = $ws.Range("Z2").Value2
$countryCa = $ws.Range("Z3").Value2
$countryAunz = $ws.Range("Z4").Value2

$ws.Range("Z1:Z4").Clear()

# Vendor / Country / Budget rows
$ws.Range("A2").Value = "Ad 4Game"
$ws.Range("B2").Value = $countryUs
$ws.Range("C2").Value = 1500

$ws.Range("A3").Value = "Ad 4Game"
$ws.Range("B3").Value = $countryAunz
$ws.Range("C3").Value = 2500

$ws.Range("A4").Value = "Ad 4Game"
$ws.Range("B4").Value = $countryBr
$ws.Range("C4").Value = 500

$ws.Range("A5").Value = "Ad 4Game"
$ws.Range("B5").Value = $countryCa
$ws.Range("C5").Value = 3500

$ws.Range("A6").Value = "Exmox"
$ws.Range("B6").Value = $countryBr
$ws.Range("C6").Value = 0

$ws.Range("A7").Value = "gdn"
$ws.Range("B7").Value = $countryAunz
$ws.Range("C7").Value = 1500

# Selection matches the saved state in the diff
$ws.Range("C6").Select()
